$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for the refreshed crypto
# snapshot. Price values that look like plain decimals are written with a
# leading apostrophe so Excel keeps them as text (matching the original
# inline-string cell type) instead of auto-converting them to numbers.
$ws.Range("D2").Value = "23.423.02"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.645.06"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'299.43"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("D7").Value = "'0.3791"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'0.3528"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").Value = "'49.93"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").Value = "'0.08078"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'22.05"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").Value = "'6.372"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "'7.314"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D17").Value = "1.648.03"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "'96.71"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "'0.06993"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "'17.36"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'12.36"
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").Value = "23.422.11"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").Value = "'2.498"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "'2.907"
$ws.Range("E26").Value = "  -5.74%  "
$ws.Range("D27").Value = "'20.85"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "'153.15"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "'5.210"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "'132.57"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "1.826.07"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "'6.875"
$ws.Range("D33").Value = "'2.120"
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").Value = "'11.41"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").Value = "'0.9807"
$ws.Range("E35").Value = "  -9.52%  "
$ws.Range("D36").Value = "'0.02702"
$ws.Range("E36").Value = "  -4.71%  "
$ws.Range("D37").Value = "'0.08752"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "'0.2429"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("D39").Value = "'5.909"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("D40").Value = "'0.06813"
$ws.Range("E40").Value = "  -4.56%  "
$ws.Range("D41").Value = "'12.84"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").Value = "'0.6859"
$ws.Range("E42").Value = "  -3.31%  "
$ws.Range("D43").Value = "'1.292"
$ws.Range("E43").Value = "  -4.10%  "
$ws.Range("D44").Value = "'15.63"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'0.6337"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'0.07716"
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").Value = "'127.02"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "'1.141"
$ws.Range("E51").Value = "  -4.67%  "
